# Update the "R Squared" column (F) values on Sheet1, rows 2-7
# to reflect the updated error measures function results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 0.9832545374323554
$ws.Range("F3").Value = 0.9978881527504846
$ws.Range("F4").Value = 0.9941867890276357
$ws.Range("F5").Value = 0.5673758213742874
$ws.Range("F6").Value = 0.8351813272024434
$ws.Range("F7").Value = 0.8146458125138485
